# Updated cryptos list on Tue Jun 27 15:54:53 UTC 2023 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) columns with the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.587.02"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.876.26"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'238.81"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.4799"
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("D8").Value = "'0.2812"
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("D9").Value = "'0.06492"
$ws.Range("D10").Value = "1.917.05"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").Value = "'0.07472"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "'16.47"
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("D13").Value = "'5.083"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").Value = "'87.82"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "'0.6620"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "30.536.27"
$ws.Range("D17").Value = "'13.23"
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'0.000007543"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("D20").Value = "2.124.80"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "'226.70"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'5.260"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("D24").Value = "'6.136"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").Value = "'9.301"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "'167.37"
$ws.Range("E26").Value = "  +2.85%  "
$ws.Range("D27").Value = "'18.45"
$ws.Range("E27").Value = "  -4.91%  "
$ws.Range("D28").Value = "'1.929"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").Value = "'1.398"
$ws.Range("E29").Value = "  -4.17%  "
$ws.Range("D30").Value = "'0.09666"
$ws.Range("E30").Value = "  +4.43%  "
$ws.Range("D31").Value = "'4.329"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "'3.993"
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").Value = "'0.05056"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").Value = "'1.215"
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("D35").Value = "'0.7465"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").Value = "'2.710"
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("D37").Value = "'0.01858"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").Value = "'0.9086"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").Value = "'2.059"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").Value = "'105.70"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "'0.4256"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").Value = "'0.9987"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").Value = "'5.749"
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("D45").Value = "'7.311"
$ws.Range("E45").Value = "  -4.57%  "
$ws.Range("D46").Value = "'0.1284"
$ws.Range("E46").Value = "  -3.59%  "
$ws.Range("D47").Value = "'63.53"
$ws.Range("E47").Value = "  -3.40%  "

# Rows 48/49 swapped rank order (EnergySwap <-> NEARProtocol) with refreshed price/volume
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.463"
$ws.Range("E48").Value = "  -8.74%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.876"
$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("D50").Value = "'33.64"
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("D51").Value = "'0.05646"
$ws.Range("E51").Value = "  -1.44%  "
